# Weekly update: a new price-report row for Cilantro (Terminal Hortofrutícola
# Agro Chillán) is inserted at row 291, pushing the existing rows 291-317
# down to 292-318 (dimension grows from R317 to R318).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 291; rows 291..317 shift to 292..318.
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new week's data.
$ws.Range("A291").Value = 7
$ws.Range("B291").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C291").Value = "Ñuble"
$ws.Range("D291").Value = 45194
$ws.Range("E291").Value = 16
$ws.Range("F291").Value = 100112040
$ws.Range("G291").Value = "Cilantro"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Primera"
$ws.Range("J291").Value = 200
$ws.Range("K291").Value = 1500
$ws.Range("L291").Value = 1500
$ws.Range("M291").Value = 1500
$ws.Range("N291").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O291").Value = "Región de Ñuble"
$ws.Range("P291").Value = 1500
$ws.Range("Q291").Value = 1
$ws.Range("R291").Value = "Hortaliza"
